# Reorder the species-record rows 9-16 (whole-row permutation, keyed by
# the original record id in column A) and then patch two AO ("microhabitat")
# descriptions that genuinely changed for the two duplicate "Bronshjon"
# records that land on rows 14 and 16 after the reorder.
#
# Only the columns whose content actually differs from row to row within
# 9-16 are touched (A,B,D,E,F,G,H,J,K,L,M,N,Q,R,AC,AF,AO); every other
# column (dates, location, observer names, booleans, ...) is identical on
# every one of these rows already, so it is left alone - this also avoids
# Excel's automatic "looks like a date" re-typing that would otherwise hit
# the Y/Z/AA/AB text cells ("2023-08-15") if we rewrote them unnecessarily.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$firstRow = 9
$lastRow  = 16

# Column indices (1-based) that vary between the rows being reordered.
$cols = @(1, 2, 4, 5, 6, 7, 8, 10, 11, 12, 13, 14, 17, 18, 29, 32, 41)

# --- 1. Snapshot the relevant cells in rows 9-16 ---------------------------
$snapshot = @{}
for ($r = $firstRow; $r -le $lastRow; $r++) {
    $rowVals = @{}
    foreach ($c in $cols) {
        $v = $ws.Cells.Item($r, $c).Value()
        if ($v -eq $null) { $v = "" }
        $rowVals[$c] = $v
    }
    $snapshot[$r] = $rowVals
}

# --- 2. Destination row -> source row mapping (whole-row permutation) -----
# Determined by matching column A (record id) between the two states.
$rowMap = @{
    9  = 15
    10 = 13
    11 = 14
    12 = 10
    13 = 16
    14 = 11
    15 = 9
    16 = 12
}

# --- 3. Write each destination row from its mapped source snapshot --------
foreach ($destRow in $rowMap.Keys) {
    $srcRow = $rowMap[$destRow]
    $vals = $snapshot[$srcRow]
    foreach ($c in $cols) {
        $ws.Cells.Item($destRow, $c).Value = $vals[$c]
    }
}

# --- 4. Fix the two microhabitat ("AO") descriptions that actually changed
#        content (not just moved) for the two "Bronshjon" records now on
#        rows 14 and 16.
$ws.Range("AO14").Value = "gammeltall"
$ws.Range("AO16").Value = "silverstubbe av tall"
